$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6,  "hello", "hello english",    "US-en English", "active"),
    @(7,  "hello", "kannada hello",    "IN-kn Kannada", "active"),
    @(8,  "bye",   "hello english",    "US-en English", "active"),
    @(9,  "bye",   "kannada hello",    "IN-kn Kannada", "active"),
    @(10, "hello", "hello in english", "US-en English", "active"),
    @(11, "bye",   "bye in english",   "US-en English", "active")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
